$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.629.74"
Set-TextValue "E2" "  +0.44%  "
Set-TextValue "D3" "1.773.09"
Set-TextValue "E3" "  +1.41%  "
Set-TextValue "D4" "0.9993"
Set-TextValue "E4" "  -0.44%  "
Set-TextValue "D5" "324.96"
Set-TextValue "E5" "  +0.18%  "
Set-TextValue "D6" "0.9971"
Set-TextValue "E6" "  -0.57%  "
Set-TextValue "D7" "0.4587"
Set-TextValue "E7" "  +3.42%  "
Set-TextValue "D8" "0.3577"
Set-TextValue "E8" "  -0.87%  "
Set-TextValue "D9" "0.07471"
Set-TextValue "E9" "  +0.07%  "
Set-TextValue "D10" "41.67"
Set-TextValue "E10" "  -1.23%  "
Set-TextValue "D11" "1.099"
Set-TextValue "E11" "  +0.27%  "
Set-TextValue "D12" "0.9987"
Set-TextValue "E12" "  -0.38%  "
Set-TextValue "D13" "20.87"
Set-TextValue "E13" "  +1.36%  "
Set-TextValue "D14" "6.024"
Set-TextValue "E14" "  -0.06%  "
Set-TextValue "D15" "7.204"
Set-TextValue "E15" "  +0.91%  "
Set-TextValue "D16" "1.767.59"
Set-TextValue "E16" "  +0.70%  "
Set-TextValue "D17" "93.77"
Set-TextValue "E17" "  +1.79%  "
Set-TextValue "D18" "0.00001056"
Set-TextValue "E18" "  -0.25%  "
Set-TextValue "D19" "0.06441"
Set-TextValue "E19" "  +0.58%  "
Set-TextValue "D20" "0.9974"
Set-TextValue "E20" "  -0.43%  "
Set-TextValue "D21" "17.11"
Set-TextValue "E21" "  +1.59%  "
Set-TextValue "D22" "5.772"
Set-TextValue "E22" "  -1.34%  "
Set-TextValue "D23" "27.734.03"
Set-TextValue "E23" "  +0.62%  "
Set-TextValue "D24" "11.26"
Set-TextValue "E24" "  +0.74%  "
Set-TextValue "D25" "2.073"
Set-TextValue "E25" "  -0.84%  "
Set-TextValue "D26" "165.27"
Set-TextValue "E26" "  +2.38%  "
Set-TextValue "D27" "20.24"
Set-TextValue "E27" "  -0.49%  "
Set-TextValue "D28" "1.967.60"
Set-TextValue "E28" "  +0.71%  "
Set-TextValue "D29" "2.166"
Set-TextValue "E29" "  +3.08%  "
Set-TextValue "D30" "126.02"
Set-TextValue "E30" "  +1.20%  "
Set-TextValue "D31" "1.090"
Set-TextValue "E31" "  +1.31%  "
Set-TextValue "D32" "0.09198"
Set-TextValue "E32" "  +2.41%  "
Set-TextValue "D33" "3.669"
Set-TextValue "E33" "  +0.58%  "
Set-TextValue "D34" "5.537"
Set-TextValue "E34" "  +0.29%  "
Set-TextValue "D35" "11.82"
Set-TextValue "E35" "  -1.14%  "
Set-TextValue "D36" "0.02286"
Set-TextValue "E36" "  -1.35%  "
Set-TextValue "D37" "0.06063"
Set-TextValue "E37" "  +1.51%  "
Set-TextValue "D38" "0.2092"
Set-TextValue "E38" "  +0.33%  "
Set-TextValue "D39" "0.6298"
Set-TextValue "E39" "  -0.54%  "
Set-TextValue "D40" "4.948"
Set-TextValue "E40" "  +0.14%  "
Set-TextValue "D41" "1.184"
Set-TextValue "E41" "  -1.78%  "
Set-TextValue "D42" "1.386"
Set-TextValue "E42" "  -0.28%  "
Set-TextValue "D43" "7.792"
Set-TextValue "E43" "  +0.41%  "
Set-TextValue "D44" "13.33"
Set-TextValue "E44" "  +0.84%  "
Set-TextValue "D45" "3.720"
Set-TextValue "E45" "  +0.16%  "
Set-TextValue "D46" "0.5882"
Set-TextValue "E46" "  +0.18%  "
Set-TextValue "D47" "122.30"
Set-TextValue "E47" "  +0.96%  "
Set-TextValue "D48" "1.946"
Set-TextValue "E48" "  +0.15%  "
Set-TextValue "D49" "0.06937"
Set-TextValue "E49" "  +1.12%  "
Set-TextValue "E50" "  -1.54%  "
Set-TextValue "D51" "72.58"
Set-TextValue "E51" "  +0.63%  "
